# Scheduled runner update: refresh market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) for the affected
# leve rows across each crafting-class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4118
$ws.Range("J51").Value = 4221.3335
$ws.Range("L51").Value = 4221.3335
$ws.Range("N51").Value = -5189.3335

$ws.Range("H74").Value = 5177.1816
$ws.Range("I74").Value = 5177.1816
$ws.Range("K74").Value = 5177.1816
$ws.Range("M74").Value = -4241.1816

$ws.Range("H77").Value = 5177.1816
$ws.Range("I77").Value = 5177.1816
$ws.Range("K77").Value = 25885.908
$ws.Range("M77").Value = -21205.908

$ws.Range("H80").Value = 1022.28
$ws.Range("I80").Value = 834.2105
$ws.Range("J80").Value = 1617.8334
$ws.Range("K80").Value = 2502.6315
$ws.Range("L80").Value = 4853.5002
$ws.Range("M80").Value = -1504.6315
$ws.Range("N80").Value = -6849.5002

$ws.Range("H83").Value = 1022.28
$ws.Range("I83").Value = 834.2105
$ws.Range("J83").Value = 1617.8334
$ws.Range("K83").Value = 7507.8945
$ws.Range("L83").Value = 14560.5006
$ws.Range("M83").Value = -2515.8945
$ws.Range("N83").Value = -24544.5006

$ws.Range("H106").Value = 77817.71000000001
$ws.Range("I106").Value = 104146.4
$ws.Range("K106").Value = 104146.4
$ws.Range("M106").Value = -103515.4

$ws.Range("H132").Value = 4349.6284
$ws.Range("I132").Value = 2459.8154
$ws.Range("K132").Value = 7379.4462
$ws.Range("M132").Value = -4849.4462

$ws.Range("H134").Value = 100195.5
$ws.Range("J134").Value = 100195.5
$ws.Range("L134").Value = 100195.5
$ws.Range("N134").Value = -110335.5

$ws.Range("H137").Value = 3594.2307
$ws.Range("I137").Value = 4627.0645
$ws.Range("J137").Value = 2069.5715
$ws.Range("K137").Value = 13881.1935
$ws.Range("L137").Value = 6208.7145
$ws.Range("M137").Value = -11331.1935
$ws.Range("N137").Value = -11308.7145

$ws.Range("H138").Value = 12988848
$ws.Range("I138").Value = 21277962
$ws.Range("J138").Value = 2569.4
$ws.Range("K138").Value = 63833886
$ws.Range("L138").Value = 7708.200000000001
$ws.Range("M138").Value = -63828746
$ws.Range("N138").Value = -17988.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4501.164
$ws.Range("I32").Value = 4151.9663
$ws.Range("J32").Value = 7076.5
$ws.Range("K32").Value = 4151.9663
$ws.Range("L32").Value = 7076.5
$ws.Range("M32").Value = -3864.9663
$ws.Range("N32").Value = -7650.5

$ws.Range("H45").Value = 7989.2
$ws.Range("I45").Value = 11172.143
$ws.Range("K45").Value = 11172.143
$ws.Range("M45").Value = -10795.143

$ws.Range("H69").Value = 180965
$ws.Range("J69").Value = 180965
$ws.Range("L69").Value = 180965
$ws.Range("N69").Value = -182463

$ws.Range("H72").Value = 180965
$ws.Range("J72").Value = 180965
$ws.Range("L72").Value = 542895
$ws.Range("N72").Value = -550383

$ws.Range("H132").Value = 1157.3125
$ws.Range("J132").Value = 1499.7142
$ws.Range("L132").Value = 4499.142599999999
$ws.Range("N132").Value = -9559.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 100554.5
$ws.Range("J42").Value = 100554.5
$ws.Range("L42").Value = 100554.5
$ws.Range("N42").Value = -101210.5

$ws.Range("H43").Value = 100554.5
$ws.Range("J43").Value = 100554.5
$ws.Range("L43").Value = 100554.5
$ws.Range("N43").Value = -100916.5

$ws.Range("H70").Value = 170966.5
$ws.Range("J70").Value = 170966.5
$ws.Range("L70").Value = 170966.5
$ws.Range("N70").Value = -171552.5

$ws.Range("H73").Value = 170966.5
$ws.Range("J73").Value = 170966.5
$ws.Range("L73").Value = 170966.5
$ws.Range("N73").Value = -172994.5

$ws.Range("H139").Value = 143213.28
$ws.Range("J139").Value = 149998.92
$ws.Range("L139").Value = 149998.92
$ws.Range("N139").Value = -160278.92

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2535.6428
$ws.Range("I31").Value = 1626
$ws.Range("J31").Value = 4809.75
$ws.Range("K31").Value = 1626
$ws.Range("L31").Value = 4809.75
$ws.Range("M31").Value = -1331
$ws.Range("N31").Value = -5399.75

$ws.Range("H34").Value = 2535.6428
$ws.Range("I34").Value = 1626
$ws.Range("J34").Value = 4809.75
$ws.Range("K34").Value = 1626
$ws.Range("L34").Value = 4809.75
$ws.Range("M34").Value = -1424
$ws.Range("N34").Value = -5213.75

$ws.Range("H134").Value = 2435.34
$ws.Range("I134").Value = 2257.4888
$ws.Range("J134").Value = 4036
$ws.Range("K134").Value = 6772.4664
$ws.Range("L134").Value = 12108
$ws.Range("M134").Value = -4237.4664
$ws.Range("N134").Value = -17178

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 827.36664
$ws.Range("J5").Value = 1397.8334
$ws.Range("L5").Value = 4193.5002
$ws.Range("N5").Value = -4417.5002

$ws.Range("H92").Value = 628.5833
$ws.Range("J92").Value = 645.8
$ws.Range("L92").Value = 1937.4
$ws.Range("N92").Value = -4433.4

$ws.Range("H98").Value = 572.2778
$ws.Range("J98").Value = 578.6429000000001
$ws.Range("L98").Value = 1735.9287
$ws.Range("N98").Value = -4731.9287

$ws.Range("H132").Value = 4273.467
$ws.Range("I132").Value = 4225.25
$ws.Range("K132").Value = 38027.25
$ws.Range("M132").Value = -35497.25

$ws.Range("H135").Value = 827.36664
$ws.Range("J135").Value = 1397.8334
$ws.Range("L135").Value = 12580.5006
$ws.Range("N135").Value = -17650.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 16979
$ws.Range("I44").Value = 18299.334
$ws.Range("J44").Value = 14998.5
$ws.Range("K44").Value = 18299.334
$ws.Range("L44").Value = 14998.5
$ws.Range("M44").Value = -17703.334
$ws.Range("N44").Value = -16190.5

$ws.Range("H99").Value = 12003.917
$ws.Range("I99").Value = 6256.5
$ws.Range("J99").Value = 23498.75
$ws.Range("K99").Value = 6256.5
$ws.Range("L99").Value = 23498.75
$ws.Range("M99").Value = -4010.5
$ws.Range("N99").Value = -27990.75

$ws.Range("H132").Value = 2625.1667
$ws.Range("I132").Value = 2381.566
$ws.Range("K132").Value = 7144.697999999999
$ws.Range("M132").Value = -4614.697999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3826.1155
$ws.Range("I16").Value = 3735.7058
$ws.Range("J16").Value = 3996.889
$ws.Range("K16").Value = 3735.7058
$ws.Range("L16").Value = 3996.889
$ws.Range("M16").Value = -3565.7058
$ws.Range("N16").Value = -4336.889

$ws.Range("H22").Value = 988.7778
$ws.Range("I22").Value = 987.25
$ws.Range("K22").Value = 987.25
$ws.Range("M22").Value = -692.25

$ws.Range("H27").Value = 988.7778
$ws.Range("I27").Value = 987.25
$ws.Range("K27").Value = 987.25
$ws.Range("M27").Value = -880.25

$ws.Range("H100").Value = 3848.3333
$ws.Range("I100").Value = 4528
$ws.Range("J100").Value = 2998.75
$ws.Range("K100").Value = 4528
$ws.Range("L100").Value = 2998.75
$ws.Range("M100").Value = -3987
$ws.Range("N100").Value = -4080.75

$ws.Range("H123").Value = 46593.75
$ws.Range("J123").Value = 46593.75
$ws.Range("L123").Value = 46593.75
$ws.Range("N123").Value = -56393.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 38000
$ws.Range("J118").Value = 38000
$ws.Range("L118").Value = 38000
$ws.Range("N118").Value = -41314
